# GDE-6676 Code and DataSet Updates
# FuncVal01: insert a new "Test_Case" column after rowid, refresh the
# Zone2/Zone3 sample data with the new Multi_E2E path + trimmed filenames,
# add the Portfolio + Business_Date columns, and make FuncVal01 the active
# tab/selection (previously FuncVal02 was active).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FuncVal01")
$ws2 = $wb.Worksheets.Item("FuncVal02")

# --- FuncVal01: insert new column B ("Test_Case") --------------------------
$ws1.Columns.Item(2).Insert()

# Restore the header style on the newly inserted B1 (Insert() copies the
# left neighbour's style, but the header row uses the shared header style).
$ws1.Range("C1").Copy()
$ws1.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row --------------------------------------------------------
$ws1.Range("B1").Value = "Test_Case"
$ws1.Range("C1").Value = "CSV_FilePath"
$ws1.Range("D1").Value = "Bal_Subledger_CSV_FileName"
$ws1.Range("E1").Value = "GL_Account_Num_CSV_FileName"
$ws1.Range("F1").Value = "GL_Short_Name_CSV_FileName"
$ws1.Range("G1").Value = "Branch_CSV_FileName"
$ws1.Range("H1").Value = "Expense_CSV_FileName"
$ws1.Range("I1").Value = "Portfolio_CSV_FileName"
$ws1.Range("J1").Value = "Business_Date"

# --- Row 2 (Zone3 sample) ------------------------------------------------
$ws1.Range("B2").Value = "DWELIQ_FuncVal01_Z3"
$ws1.Range("C2").Value = "C:\Git_Evergreen\fms_cba\DataSet\Integration_DataSet\Extracts\DWE_LIQ\DWE_LIQ_Extracts\DWELIQ_Multi_E2E_001\ZONE3\CCB_LIQ_SYD_"
$ws1.Range("D2").Value = "CCB_LIQ_SYD_VLS_BAL_SUBLEDGER_"
$ws1.Range("E2").Value = "CCB_LIQ_SYD_VLS_GL_ACCT_NUM_"
$ws1.Range("F2").Value = "CCB_LIQ_SYD_VLS_GL_SHORT_NAME_"
$ws1.Range("G2").Value = "CCB_LIQ_SYD_VLS_BRANCH_"
$ws1.Range("H2").Value = "CCB_LIQ_SYD_VLS_EXPENSE_"
$ws1.Range("I2").Value = "CCB_LIQ_SYD_VLS_PORTFOLIO_"
# Business_Date must stay a text value (not get coerced to a number) and
# keep the default (unstyled) cell format, so copy the identical text
# value that FuncVal02!F2 already holds rather than typing a numeric
# literal into J2 directly.
$ws2.Range("F2").Copy()
$ws1.Range("J2").PasteSpecial(-4163)

# --- Row 3 (Zone2 sample) ------------------------------------------------
$ws1.Range("B3").Value = "DWELIQ_FuncVal01_Z2"
$ws1.Range("C3").Value = "C:\Git_Evergreen\fms_cba\DataSet\Integration_DataSet\Extracts\DWE_LIQ\DWE_LIQ_Extracts\DWELIQ_Multi_E2E_002\ZONE2\CCB_LIQ_EUR_"
$ws1.Range("D3").Value = "CCB_LIQ_EUR_VLS_BAL_SUBLEDGER_"
$ws1.Range("E3").Value = "CCB_LIQ_EUR_VLS_GL_ACCT_NUM_"
$ws1.Range("F3").Value = "CCB_LIQ_EUR_VLS_GL_SHORT_NAME_"
$ws1.Range("G3").Value = "CCB_LIQ_EUR_VLS_BRANCH_"
$ws1.Range("H3").Value = "CCB_LIQ_EUR_VLS_EXPENSE_"
$ws1.Range("I3").Value = "CCB_LIQ_EUR_VLS_PORTFOLIO_"
$ws2.Range("F3").Copy()
$ws1.Range("J3").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Column widths (best effort; matches the author's resize of the new
# Test_Case/CSV_FilePath columns) ---------------------------------------
$ws1.Columns.Item(2).ColumnWidth = 21.42578125
$ws1.Columns.Item(3).ColumnWidth = 129.85546875

# --- Sheet selection / active tab ---------------------------------------
# FuncVal02 loses its "last selected" cell/tab; FuncVal01 becomes active.
$ws2.Range("F1").Select()
$ws1.Activate()
$ws1.Range("D7").Select()
